$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$c = $ws.Range("D2")
$c.Value = '''43.220.17'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -2.19%  '

# Row 3 - Ethereum
$c = $ws.Range("D3")
$c.Value = '''2.334.49'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +2.85%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  -0.04%  '

# Row 5 - XRP
$c = $ws.Range("D5")
$c.Value = '''0.648'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.45%  '

# Row 6 - BNB
$c = $ws.Range("D6")
$c.Value = '''230.85'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.10%  '

# Row 7 - Solana
$c = $ws.Range("D7")
$c.Value = '''64.98'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.94%  '

# Row 8 - USDC
$ws.Range("E8").Value = '  -0.04%  '

# Row 9 - Cardano
$c = $ws.Range("D9")
$c.Value = '''0.452'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.46%  '

# Row 10 - Dogecoin
$c = $ws.Range("D10")
$c.Value = '''0.0948'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -5.08%  '

# Row 11 - OKB
$c = $ws.Range("D11")
$c.Value = '''56.72'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.89%  '

# Row 12 - Avalanche
$c = $ws.Range("D12")
$c.Value = '''26.57'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -3.54%  '

# Row 13 - WrappedliquidstakedEther2.0
$c = $ws.Range("D13")
$c.Value = '''2.680.01'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +2.81%  '

# Row 14 - TRON
$ws.Range("E14").Value = '  -1.49%  '

# Row 15 - Chainlink
$c = $ws.Range("D15")
$c.Value = '''15.24'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -3.70%  '

# Row 16 - Polkadot
$ws.Range("E16").Value = '  +1.71%  '

# Row 17 - Polygon
$c = $ws.Range("D17")
$c.Value = '''0.836'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.61%  '

# Row 18 - WrappedEther
$c = $ws.Range("D18")
$c.Value = '''2.335.68'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +2.36%  '

# Row 19 - WrappedBTC
$c = $ws.Range("D19")
$c.Value = '''43.138.66'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -2.01%  '

# Row 20 - ShibaInu
$c = $ws.Range("D20")
$c.Value = "'0.0₃0970"
$c.Style = "Normal"
$ws.Range("E20").Value = '  -4.12%  '

# Row 21 - Litecoin
$c = $ws.Range("D21")
$c.Value = '''73.50'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.53%  '

# Row 22 - Uniswap
$c = $ws.Range("D22")
$c.Value = '''6.15'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.52%  '

# Row 23 - BitcoinCash
$c = $ws.Range("D23")
$c.Value = '''247.33'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -2.41%  '

# Row 24 - WEMIXToken
$ws.Range("E24").Value = '  +19.91%  '

# Row 25 - Dai
$ws.Range("E25").Value = '  +0.04%  '

# Row 26 - PancakeSwap
$c = $ws.Range("D26")
$c.Value = '''2.42'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.69%  '

# Row 27 - Toncoin
$c = $ws.Range("D27")
$c.Value = '''2.27'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.93%  '

# Row 28 - Cosmos
$c = $ws.Range("D28")
$c.Value = '''9.81'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -3.33%  '

# Row 29 - Monero
$c = $ws.Range("D29")
$c.Value = '''174.52'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +1.72%  '

# Row 30 - EthereumClassic
$c = $ws.Range("D30")
$c.Value = '''22.11'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +5.53%  '

# Row 31 - ImmutableX
$ws.Range("E31").Value = '  +4.44%  '

# Row 32 - Kaspa
$c = $ws.Range("D32")
$c.Value = '''0.128'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -8.73%  '

# Row 33 - Stellar
$ws.Range("E33").Value = '  -0.11%  '

# Row 34 - Filecoin
$c = $ws.Range("D34")
$c.Value = '''4.99'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +3.66%  '

# Row 35 - Hedera
$ws.Range("E35").Value = '  -3.58%  '

# Row 36 - InternetComputer(DFINITY)
$c = $ws.Range("D36")
$c.Value = '''4.95'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +0.97%  '

# Row 37 - LidoDAOToken
$c = $ws.Range("D37")
$c.Value = '''2.46'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +6.04%  '

# Row 38 - THORChain
$ws.Range("E38").Value = '  -1.43%  '

# Row 39 - RenderToken
$c = $ws.Range("D39")
$c.Value = '''3.56'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -6.73%  '

# Row 40 - VeChain
$ws.Range("E40").Value = '  -4.38%  '

# Row 41 - BinanceUSD
$ws.Range("E41").Value = '  +0.12%  '

# Row 42 - FraxShare
$c = $ws.Range("D42")
$c.Value = '''8.84'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +7.25%  '

# Row 43 - InjectiveProtocol
$c = $ws.Range("D43")
$c.Value = '''17.77'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +1.41%  '

# Row 44 - ARBITRUM
$ws.Range("E44").Value = '  +5.89%  '

# Row 45 - Aave
$c = $ws.Range("D45")
$c.Value = '''98.11'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.40%  '

# Row 46 - TrustWalletToken
$ws.Range("E46").Value = '  -1.51%  '

# Row 47 - FTXToken
$c = $ws.Range("D47")
$c.Value = '''4.37'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.16%  '

# Row 48 - Cronos
$c = $ws.Range("D48")
$c.Value = '''0.0940'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -5.13%  '

# Row 49 - Maker
$c = $ws.Range("D49")
$c.Value = '''1.433.26'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.00%  '

# Row 50 - Celestia
$c = $ws.Range("D50")
$c.Value = '''9.80'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -6.02%  '

# Row 51 - TerraClassic -> RocketPoolETH
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$c = $ws.Range("D51")
$c.Value = '''2.552.40'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +2.81%  '
